$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3149.1667
$ws.Range("I64").Value = 2399.2856
$ws.Range("J64").Value = 4199
$ws.Range("K64").Value = 2399.2856
$ws.Range("L64").Value = 4199
$ws.Range("M64").Value = -2151.2856
$ws.Range("N64").Value = -4695

# Row 67
$ws.Range("H67").Value = 3149.1667
$ws.Range("I67").Value = 2399.2856
$ws.Range("J67").Value = 4199
$ws.Range("K67").Value = 2399.2856
$ws.Range("L67").Value = 4199
$ws.Range("M67").Value = -1541.2856
$ws.Range("N67").Value = -5915

# Row 98
$ws.Range("H98").Value = 39178.75
$ws.Range("I98").Value = 46407.4
$ws.Range("J98").Value = 3035.5
$ws.Range("K98").Value = 46407.4
$ws.Range("L98").Value = 3035.5
$ws.Range("M98").Value = -44909.4
$ws.Range("N98").Value = -6031.5

# Row 107
$ws.Range("H107").Value = 1868
$ws.Range("J107").Value = 2545.875
$ws.Range("L107").Value = 2545.875
$ws.Range("N107").Value = -6385.875

# Row 111
$ws.Range("H111").Value = 1027.4286
$ws.Range("I111").Value = 1035
$ws.Range("J111").Value = 999.6667
$ws.Range("K111").Value = 3105
$ws.Range("L111").Value = 2999.0001
$ws.Range("M111").Value = -38
$ws.Range("N111").Value = -9133.000100000001

# Row 122
$ws.Range("H122").Value = 39178.75
$ws.Range("I122").Value = 46407.4
$ws.Range("J122").Value = 3035.5
$ws.Range("K122").Value = 139222.2
$ws.Range("L122").Value = 9106.5
$ws.Range("M122").Value = -136772.2
$ws.Range("N122").Value = -14006.5

# Row 133
$ws.Range("H133").Value = 165110.6
$ws.Range("I133").Value = 80000
$ws.Range("J133").Value = 186388.25
$ws.Range("K133").Value = 80000
$ws.Range("L133").Value = 186388.25
$ws.Range("M133").Value = -74940
$ws.Range("N133").Value = -196508.25

# Row 137
$ws.Range("H137").Value = 178202.89
$ws.Range("I137").Value = 2243.2
$ws.Range("J137").Value = 373713.66
$ws.Range("K137").Value = 6729.599999999999
$ws.Range("L137").Value = 1121140.98
$ws.Range("M137").Value = -4179.599999999999
$ws.Range("N137").Value = -1126240.98

# Row 138
$ws.Range("H138").Value = 2764.353
$ws.Range("I138").Value = 2385.5386
$ws.Range("J138").Value = 3995.5
$ws.Range("K138").Value = 7156.6158
$ws.Range("L138").Value = 11986.5
$ws.Range("M138").Value = -2016.6158
$ws.Range("N138").Value = -22266.5

# Row 139
$ws.Range("H139").Value = 62150.42
$ws.Range("I139").Value = 42499.5
$ws.Range("J139").Value = 64462.293
$ws.Range("K139").Value = 42499.5
$ws.Range("L139").Value = 64462.293
$ws.Range("M139").Value = -37359.5
$ws.Range("N139").Value = -74742.29300000001

# Row 140
$ws.Range("H140").Value = 91407.25
$ws.Range("J140").Value = 77553.37
$ws.Range("L140").Value = 77553.37
$ws.Range("N140").Value = -87913.37

# Row 141
$ws.Range("H141").Value = 9948.25
$ws.Range("I141").Value = 888
$ws.Range("J141").Value = 12968.333
$ws.Range("K141").Value = 2664
$ws.Range("L141").Value = 38904.999
$ws.Range("M141").Value = 2516
$ws.Range("N141").Value = -49264.999

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 2541.768
$ws.Range("I132").Value = 2412.9268
$ws.Range("K132").Value = 7238.780400000001
$ws.Range("M132").Value = -4708.780400000001

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 16939.904
$ws.Range("I107").Value = 17636.9
$ws.Range("K107").Value = 17636.9
$ws.Range("M107").Value = -15716.9

# Row 134
$ws.Range("H134").Value = 2501.0454
$ws.Range("I134").Value = 2373.6296
$ws.Range("K134").Value = 7120.888800000001
$ws.Range("M134").Value = -4585.888800000001

# Row 138
$ws.Range("H138").Value = 97939.8
$ws.Range("J138").Value = 99924.75
$ws.Range("L138").Value = 99924.75
$ws.Range("N138").Value = -110204.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2101.57
$ws.Range("I31").Value = 1745.1132
$ws.Range("J31").Value = 2503.532
$ws.Range("K31").Value = 1745.1132
$ws.Range("L31").Value = 2503.532
$ws.Range("M31").Value = -1450.1132
$ws.Range("N31").Value = -3093.532

# Row 34
$ws.Range("H34").Value = 2101.57
$ws.Range("I34").Value = 1745.1132
$ws.Range("J34").Value = 2503.532
$ws.Range("K34").Value = 1745.1132
$ws.Range("L34").Value = 2503.532
$ws.Range("M34").Value = -1543.1132
$ws.Range("N34").Value = -2907.532

# Row 58
$ws.Range("H58").Value = 3148
$ws.Range("I58").Value = 2835.762
$ws.Range("K58").Value = 2835.762
$ws.Range("M58").Value = -2632.762

# Row 62
$ws.Range("H62").Value = 134137.12
$ws.Range("I62").Value = 253374.75
$ws.Range("J62").Value = 14899.5
$ws.Range("K62").Value = 253374.75
$ws.Range("L62").Value = 14899.5
$ws.Range("M62").Value = -252750.75
$ws.Range("N62").Value = -16147.5

# Row 65
$ws.Range("H65").Value = 134137.12
$ws.Range("I65").Value = 253374.75
$ws.Range("J65").Value = 14899.5
$ws.Range("K65").Value = 1266873.75
$ws.Range("L65").Value = 74497.5
$ws.Range("M65").Value = -1263753.75
$ws.Range("N65").Value = -80737.5

# Row 132
$ws.Range("H132").Value = 4762
$ws.Range("I132").Value = 2306.4119
$ws.Range("J132").Value = 46507
$ws.Range("K132").Value = 6919.2357
$ws.Range("L132").Value = 139521
$ws.Range("M132").Value = -4389.2357
$ws.Range("N132").Value = -144581

# Row 134
$ws.Range("H134").Value = 3271.484
$ws.Range("I134").Value = 3197.963
$ws.Range("J134").Value = 3767.75
$ws.Range("K134").Value = 9593.889000000001
$ws.Range("L134").Value = 11303.25
$ws.Range("M134").Value = -7058.889000000001
$ws.Range("N134").Value = -16373.25

# Row 135
$ws.Range("H135").Value = 67980
$ws.Range("J135").Value = 67980
$ws.Range("L135").Value = 67980
$ws.Range("M135").Value = -78120

# Row 136
$ws.Range("H136").Value = 3148
$ws.Range("I136").Value = 2835.762
$ws.Range("K136").Value = 8507.286
$ws.Range("M136").Value = -5957.286

# Row 138
$ws.Range("H138").Value = 81595.89
$ws.Range("J138").Value = 81595.89
$ws.Range("L138").Value = 81595.89
$ws.Range("N138").Value = -91875.89

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 65755
$ws.Range("I102").Value = 116322.664
$ws.Range("J102").Value = 15187.333
$ws.Range("K102").Value = 116322.664
$ws.Range("L102").Value = 15187.333
$ws.Range("M102").Value = -114700.664
$ws.Range("N102").Value = -18431.333

# Row 135
$ws.Range("H135").Value = 99364
$ws.Range("J135").Value = 99364
$ws.Range("L135").Value = 99364
$ws.Range("N135").Value = -109504

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 761.46155
$ws.Range("I107").Value = 721.5
$ws.Range("J107").Value = 894.6667
$ws.Range("K107").Value = 2164.5
$ws.Range("L107").Value = 2684.0001
$ws.Range("M107").Value = -244.5
$ws.Range("N107").Value = -6524.0001

# Row 113
$ws.Range("H113").Value = 817.5714
$ws.Range("I113").Value = 856.4
$ws.Range("J113").Value = 494
$ws.Range("K113").Value = 2569.2
$ws.Range("L113").Value = 1482
$ws.Range("M113").Value = -399.1999999999998
$ws.Range("N113").Value = -5822

# Row 139
$ws.Range("H139").Value = 84815.89
$ws.Range("J139").Value = 84815.89
$ws.Range("L139").Value = 84815.89
$ws.Range("N139").Value = -95095.89

# Row 141
$ws.Range("H141").Value = 121744.25
$ws.Range("J141").Value = 121744.25
$ws.Range("L141").Value = 121744.25
$ws.Range("N141").Value = -132104.25

Write-Host "Applied all Phoenix_Profits updates"
